$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Append the new sentences about the ASP.Net implementation right after
#    "... den Zeitrechner entschieden." as a brand-new run in the same
#    paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("den Zeitrechner entschieden.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter(" Es soll in ASP.Net umgesetzt werden. Es gibt eine Webseite und eine API. Die Webseite wird ein Datum an die API senden und die wird die Differenz von der Zeit zu der aktuellen Zeit ausrechnen. Das Ergebnis wird dann zurückgegeben und von der Webseite angezeigt.")
}

# ---------------------------------------------------------------------------
# 2) In the work-log table (the 2nd table in the document) the page
#    repaginates after the text above grew, so Word's rendered-page-break
#    marker slides from the "Backend/120" cells up onto the "2" cell.
#    Row 3 is the "2" / "2.A" / "19.1.24" / "Liam" / "Man erstellt ein
#    Backend..." / "120" row.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(2)
$row = 3

# --- Column 1 ("2"): add <w:lastRenderedPageBreak/> before the text run ----
$cellNr = $t.Cell($row, 1)
$paraNr = $cellNr.Range.Paragraphs.Item(1)
$xmlNr = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="095AAC32" w14:textId="43D09DC9" w:rsidR="0045407F" w:rsidRDefault="001E0AEF"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>2</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$paraNr.Range.InsertXML($xmlNr)

# --- Column 5 ("Man erstellt ein Backend..."): merge the two runs into one
#     run and drop the now-stale <w:lastRenderedPageBreak/> -----------------
$cellDesc = $t.Cell($row, 5)
$paraDesc = $cellDesc.Range.Paragraphs.Item(1)
$xmlDesc = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="25C78362" w14:textId="73C62E0B" w:rsidR="0045407F" w:rsidRDefault="00B368E0"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>Man erstellt ein Backend, dass durch das Programm aufgerufen werden kann und eine Zeit &#252;bergeben kann.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$paraDesc.Range.InsertXML($xmlDesc)

# --- Column 6 ("120"): remove the stale <w:lastRenderedPageBreak/> --------
$cell120 = $t.Cell($row, 6)
$para120 = $cell120.Range.Paragraphs.Item(1)
$xml120 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7207FE92" w14:textId="5D8C4262" w:rsidR="0045407F" w:rsidRDefault="00B368E0"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t>120</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para120.Range.InsertXML($xml120)

Write-Output "done"
